$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number (column C = plain_password) to new value
$passwords = @{
    2  = "Pmgt24"
    3  = "ZWMx89"
    4  = "7DvZ84"
    5  = "Y6Pw28"
    6  = "xLEF54"
    7  = "nIEq25"
    8  = "Vbbz45"
    9  = "qQSl25"
    10 = "17PB68"
    11 = "D8VV79"
    12 = "ZLOU82"
    13 = "DE8U13"
    14 = "dRRw50"
    15 = "IySy70"
    16 = "XnNF94"
    17 = "FAgr66"
    18 = "pyJ358"
    19 = "tLks33"
    20 = "h7Ls63"
    21 = "0Plf66"
    22 = "f1NL68"
    23 = "17fc48"
}

foreach ($row in $passwords.Keys) {
    $ws.Cells.Item($row, 3).Value = $passwords[$row]
}
